# Excel Reader and Data driven test has been added
# - Refreshes the Login sheet's sample credentials
# - Adds a "ManageListing" sheet (Title / Deleteaction sample row)
# - Adds a "ShareSkill" sheet (Title / Description / Tag / SKTag sample row,
#   with spare date/time formatted helper columns)

$wb = $excel.ActiveWorkbook
$white = 16777215

# --- Update existing "Login" sheet -----------------------------------------
$loginWs = $wb.Worksheets.Item("Login")
$loginWs.Range("A2").Value = "451075672@qq.com"
$loginWs.Range("B2").Value = "abc123"

$loginWs.Range("A1:B1").Interior.Color = $white
$loginWs.Range("A1:B1").Borders.LineStyle = 1
$loginWs.Range("A2:B2").Borders.LineStyle = 1
$loginWs.Range("A1").Select()

# --- Add "ManageListing" sheet ----------------------------------------------
$manageWs = $wb.Worksheets.Add()
$manageWs.Name = "ManageListing"
$manageWs.Range("A1").Value = "Title"
$manageWs.Range("B1").Value = "Deleteaction"
$manageWs.Range("A2").Value = "Jazz"
$manageWs.Range("B2").Value = "Yes"

$manageWs.Range("A1:B1").Interior.Color = $white
$manageWs.Range("A1:B1").Borders.LineStyle = 1
$manageWs.Range("A2:B2").Borders.LineStyle = 1

$manageWs.Columns.Item(2).ColumnWidth = 14
$manageWs.Range("F16").Select()
$manageWs.Move($null, $wb.Worksheets.Item($wb.Worksheets.Count))

# --- Add "ShareSkill" sheet --------------------------------------------------
$shareWs = $wb.Worksheets.Add()
$shareWs.Name = "ShareSkill"
$shareWs.Range("A1").Value = "Title"
$shareWs.Range("B1").Value = "Description"
$shareWs.Range("C1").Value = "Tag"
$shareWs.Range("D1").Value = "SKTag"
$shareWs.Range("A2").Value = "Jazz Club"
$shareWs.Range("B2").Value = "We are true Jazz lover, If you are intrested in Jazz please feel free to join us!!"
$shareWs.Range("C2").Value = "Jazz"
$shareWs.Range("D2").Value = "Photo Skill"

$shareWs.Range("A1:D1").Interior.Color = $white
$shareWs.Range("A1:D1").Borders.LineStyle = 1
$shareWs.Range("E1:M1").Interior.Color = $white
$shareWs.Range("A2:D2").Borders.LineStyle = 1

$shareWs.Range("E2:G2").NumberFormat = "dd/mm/yyyy;@"
$shareWs.Range("H2:I2").NumberFormat = "hh:mm:ss;@"
$shareWs.Range("J2:M2").NumberFormat = "General"

$shareWs.Columns.Item(1).ColumnWidth = 10.6363636363636
$shareWs.Columns.Item(2).ColumnWidth = 92.8181818181818
$shareWs.Columns.Item(3).ColumnWidth = 17.3636363636364
$shareWs.Columns.Item(4).ColumnWidth = 14
$shareWs.Columns.Item(5).ColumnWidth = 11.8181818181818
$shareWs.Columns.Item(6).ColumnWidth = 11.8181818181818
$shareWs.Columns.Item(7).ColumnWidth = 10.6363636363636
$shareWs.Columns.Item(8).ColumnWidth = 10.6363636363636
$shareWs.Columns.Item(10).ColumnWidth = 16.2727272727273
$shareWs.Columns.Item(11).ColumnWidth = 21.9090909090909
$shareWs.Columns.Item(12).ColumnWidth = 7.54545454545455

$shareWs.Range("C22").Select()
$shareWs.Move($null, $wb.Worksheets.Item($wb.Worksheets.Count))

$shareWs.Activate()
$excel.ActiveWindow.WindowState = -4143
